# Minor dataset queries cleanup
# Normalize header row casing (mostly lower-casing, keep "Reorder Level" as-is)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "quantity in stock"
$ws.Range("E1").Value = "shelf Location"
$ws.Range("F1").Value = "Reorder Level"
$ws.Range("G1").Value = "Last restock date"
$ws.Range("H1").Value = "unit price"
$ws.Range("B1").Value = "Part description"
$ws.Range("C1").Value = "Supplier information"
$ws.Range("A1").Value = "Part number"
